$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as per upstream refresh.
# Column D (Price) values are forced to text to preserve exact formatting
# (e.g. "1.00", "69.197.19") the same way the source data stores them as inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.197.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.670.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "671.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.50%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.73%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.289.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.670.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.152.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.816.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.65%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.86%  "
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  -5.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.664.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0895"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.73%  "
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.15%  "
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.01%  "
$ws.Range("E49").Value = "  -6.09%  "
$ws.Range("E50").Value = "  -3.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.53%  "
